$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TE1 dup-detector test data: change the FRA reporting period for row 13
# from 202504 to 202505 (adds a new SSN-all-9s dup test case with a
# distinct period so it doesn't collide with the other 202504 rows).
$ws.Range("A13").Value = 202505

# Column B was resized (the SSN column) while editing the sheet.
$ws.Columns.Item(2).ColumnWidth = 18.75

# Leave the cursor on A14, matching where editing left off.
$ws.Range("A14").Select()
